# Updated cryptos list (prices & 1h volume change) per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.566.72"
$ws.Range("E2").Value = "'  +5.44%  "
$ws.Range("D3").Value = "'1.726.25"
$ws.Range("E3").Value = "'  +4.31%  "
$ws.Range("E4").Value = "'  +0.12%  "
$ws.Range("D5").Value = "'225.92"
$ws.Range("E5").Value = "'  +3.26%  "
$ws.Range("D6").Value = "'0.5391"
$ws.Range("E6").Value = "'  +2.91%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "'  +0.05%  "
$ws.Range("D8").Value = "'0.2682"
$ws.Range("E8").Value = "'  +0.86%  "
$ws.Range("D9").Value = "'0.06622"
$ws.Range("E9").Value = "'  +4.15%  "
$ws.Range("D10").Value = "'21.83"
$ws.Range("E10").Value = "'  +6.48%  "
$ws.Range("D11").Value = "'0.07731"
$ws.Range("E11").Value = "'  +0.52%  "
$ws.Range("D12").Value = "'4.623"
$ws.Range("E12").Value = "'  -0.34%  "
$ws.Range("D13").Value = "'1.722.78"
$ws.Range("E13").Value = "'  +5.81%  "
$ws.Range("D14").Value = "'1.962.43"
$ws.Range("E14").Value = "'  +4.27%  "
$ws.Range("D15").Value = "'0.5867"
$ws.Range("E15").Value = "'  +4.50%  "
$ws.Range("D16").Value = "'0.0₅8320"
$ws.Range("E16").Value = "'  +1.49%  "
$ws.Range("D17").Value = "'68.10"
$ws.Range("E17").Value = "'  +3.89%  "
$ws.Range("D18").Value = "'27.578.90"
$ws.Range("D19").Value = "'221.82"
$ws.Range("E19").Value = "'  +15.27%  "
$ws.Range("E20").Value = "'  +0.07%  "
$ws.Range("D21").Value = "'4.746"
$ws.Range("E21").Value = "'  +1.95%  "
$ws.Range("D22").Value = "'10.70"
$ws.Range("E22").Value = "'  +2.10%  "
$ws.Range("D23").Value = "'6.108"
$ws.Range("E23").Value = "'  +2.48%  "
$ws.Range("D24").Value = "'1.004"
$ws.Range("E24").Value = "'  +0.13%  "
$ws.Range("D25").Value = "'148.28"
$ws.Range("E25").Value = "'  +2.04%  "
$ws.Range("D26").Value = "'0.1238"
$ws.Range("E26").Value = "'  +3.50%  "
$ws.Range("D27").Value = "'1.687"
$ws.Range("E27").Value = "'  +11.31%  "
$ws.Range("D28").Value = "'7.417"
$ws.Range("E28").Value = "'  +2.10%  "
$ws.Range("D29").Value = "'16.70"
$ws.Range("E29").Value = "'  +4.54%  "
$ws.Range("D30").Value = "'0.05584"
$ws.Range("E30").Value = "'  +1.96%  "
$ws.Range("D31").Value = "'1.304"
$ws.Range("E31").Value = "'  +2.62%  "
$ws.Range("D32").Value = "'3.554"
$ws.Range("E32").Value = "'  +2.63%  "
$ws.Range("D33").Value = "'3.459"
$ws.Range("E33").Value = "'  +2.38%  "
$ws.Range("D34").Value = "'1.663"
$ws.Range("E34").Value = "'  +6.62%  "
$ws.Range("D35").Value = "'0.9660"
$ws.Range("E35").Value = "'  +1.21%  "
$ws.Range("D36").Value = "'2.828"
$ws.Range("E36").Value = "'  +1.70%  "
$ws.Range("E37").Value = "'  +1.82%  "
$ws.Range("D38").Value = "'0.5961"
$ws.Range("E38").Value = "'  +5.28%  "
$ws.Range("D39").Value = "'0.01649"
$ws.Range("E39").Value = "'  +4.28%  "
$ws.Range("D40").Value = "'5.932"
$ws.Range("E40").Value = "'  +1.08%  "
$ws.Range("D41").Value = "'0.8568"
$ws.Range("E41").Value = "'  +2.46%  "
$ws.Range("D42").Value = "'1.056.80"
$ws.Range("E42").Value = "'  +2.69%  "
$ws.Range("E43").Value = "'  +0.09%  "
$ws.Range("D44").Value = "'101.56"
$ws.Range("E44").Value = "'  +0.28%  "
$ws.Range("D45").Value = "'1.868.09"
$ws.Range("D46").Value = "'0.0₈115"
$ws.Range("E46").Value = "'  +5.03%  "
$ws.Range("D47").Value = "'59.23"
$ws.Range("E47").Value = "'  +2.33%  "
$ws.Range("D48").Value = "'8.196"
$ws.Range("E48").Value = "'  +2.59%  "
$ws.Range("D49").Value = "'0.4441"
$ws.Range("E49").Value = "'  +2.34%  "
$ws.Range("D50").Value = "'1.003"
$ws.Range("E50").Value = "'  +0.33%  "
$ws.Range("D51").Value = "'0.05275"
$ws.Range("E51").Value = "'  +1.58%  "
